$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 32666.334
$ws.Range("J108").Value = 32666.334
$ws.Range("L108").Value = 32666.334
$ws.Range("N108").Value = -40346.334
$ws.Range("H132").Value = 3965.1304
$ws.Range("I132").Value = 4122.1113
$ws.Range("J132").Value = 3400
$ws.Range("K132").Value = 12366.3339
$ws.Range("L132").Value = 10200
$ws.Range("M132").Value = -9836.333899999998
$ws.Range("N132").Value = -15260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3051.7727
$ws.Range("J45").Value = 3513.652
$ws.Range("L45").Value = 3513.652
$ws.Range("N45").Value = -4267.652
$ws.Range("H61").Value = 1537.2593
$ws.Range("I61").Value = 729.79486
$ws.Range("K61").Value = 729.79486
$ws.Range("M61").Value = -517.79486
$ws.Range("H136").Value = 1537.2593
$ws.Range("I136").Value = 729.79486
$ws.Range("K136").Value = 2189.38458
$ws.Range("M136").Value = 360.6154200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 20050
$ws.Range("J44").Value = 20050
$ws.Range("L44").Value = 20050
$ws.Range("N44").Value = -21044
$ws.Range("H99").Value = 1705.3077
$ws.Range("I99").Value = 1706.2727
$ws.Range("J99").Value = 1700
$ws.Range("K99").Value = 1706.2727
$ws.Range("L99").Value = 1700
$ws.Range("M99").Value = -208.2727
$ws.Range("N99").Value = -4696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 702.1667
$ws.Range("I107").Value = 550
$ws.Range("J107").Value = 1006.5
$ws.Range("K107").Value = 550
$ws.Range("L107").Value = 1006.5
$ws.Range("M107").Value = 1370
$ws.Range("N107").Value = -4846.5
$ws.Range("H134").Value = 1146.0869
$ws.Range("I134").Value = 842.8
$ws.Range("J134").Value = 1714.75
$ws.Range("K134").Value = 2528.4
$ws.Range("L134").Value = 5144.25
$ws.Range("M134").Value = 6.600000000000364
$ws.Range("N134").Value = -10214.25
$ws.Range("H135").Value = 50630
$ws.Range("J135").Value = 50630
$ws.Range("L135").Value = 50630
$ws.Range("N135").Value = -60770

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1322.5555
$ws.Range("J5").Value = 1820
$ws.Range("L5").Value = 5460
$ws.Range("N5").Value = -5684
$ws.Range("H80").Value = 2820
$ws.Range("J80").Value = 2820
$ws.Range("L80").Value = 8460
$ws.Range("N80").Value = -10332
$ws.Range("H83").Value = 2820
$ws.Range("J83").Value = 2820
$ws.Range("L83").Value = 25380
$ws.Range("N83").Value = -34740
$ws.Range("H121").Value = 17202.334
$ws.Range("J121").Value = 33899.668
$ws.Range("L121").Value = 101699.004
$ws.Range("N121").Value = -104319.004
$ws.Range("H122").Value = 423.14285
$ws.Range("I122").Value = 423.14285
$ws.Range("K122").Value = 3808.28565
$ws.Range("M122").Value = -1358.28565
$ws.Range("H131").Value = 754.95
$ws.Range("I131").Value = 388.33334
$ws.Range("J131").Value = 778.3511
$ws.Range("K131").Value = 1165.00002
$ws.Range("L131").Value = 2335.0533
$ws.Range("M131").Value = 3874.99998
$ws.Range("N131").Value = -12415.0533
$ws.Range("H135").Value = 1322.5555
$ws.Range("J135").Value = 1820
$ws.Range("L135").Value = 16380
$ws.Range("N135").Value = -21450
$ws.Range("H139").Value = 1724.7391
$ws.Range("I139").Value = 1149.5
$ws.Range("J139").Value = 3795.6
$ws.Range("K139").Value = 3448.5
$ws.Range("L139").Value = 11386.8
$ws.Range("M139").Value = 1691.5
$ws.Range("N139").Value = -21666.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4255.154
$ws.Range("I97").Value = 1729.4445
$ws.Range("K97").Value = 1729.4445
$ws.Range("M97").Value = -1233.4445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3148.087
$ws.Range("I7").Value = 3493.3572
$ws.Range("J7").Value = 2611
$ws.Range("K7").Value = 3493.3572
$ws.Range("L7").Value = 2611
$ws.Range("M7").Value = -3381.3572
$ws.Range("N7").Value = -2835
$ws.Range("H22").Value = 3867
$ws.Range("I22").Value = 10001
$ws.Range("K22").Value = 10001
$ws.Range("M22").Value = -9706
$ws.Range("H27").Value = 3867
$ws.Range("I27").Value = 10001
$ws.Range("K27").Value = 10001
$ws.Range("M27").Value = -9894
$ws.Range("H40").Value = 12001.667
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 12001.667
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 12001.667
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -12273.667
$ws.Range("H46").Value = 2640
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 2850
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 2850
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -3226
$ws.Range("H55").Value = 221
$ws.Range("I55").Value = 170
$ws.Range("J55").Value = 241.4
$ws.Range("K55").Value = 170
$ws.Range("L55").Value = 241.4
$ws.Range("M55").Value = 3
$ws.Range("N55").Value = -587.4
$ws.Range("H61").Value = 3532.1177
$ws.Range("I61").Value = 1588.8334
$ws.Range("J61").Value = 8196
$ws.Range("K61").Value = 1588.8334
$ws.Range("L61").Value = 8196
$ws.Range("M61").Value = -1386.8334
$ws.Range("N61").Value = -8600
$ws.Range("H93").Value = 3145.1667
$ws.Range("I93").Value = 3085.6365
$ws.Range("K93").Value = 3085.6365
$ws.Range("M93").Value = -1837.6365
$ws.Range("H113").Value = 3532.1177
$ws.Range("I113").Value = 1588.8334
$ws.Range("J113").Value = 8196
$ws.Range("K113").Value = 1588.8334
$ws.Range("L113").Value = 8196
$ws.Range("M113").Value = 581.1666
$ws.Range("N113").Value = -12536
$ws.Range("H122").Value = 1157094.6
$ws.Range("I122").Value = 3924192.8
$ws.Range("J122").Value = 4137
$ws.Range("K122").Value = 11772578.4
$ws.Range("L122").Value = 12411
$ws.Range("M122").Value = -11770128.4
$ws.Range("N122").Value = -17311
$ws.Range("H126").Value = 3148.087
$ws.Range("I126").Value = 3493.3572
$ws.Range("J126").Value = 2611
$ws.Range("K126").Value = 10480.0716
$ws.Range("L126").Value = 7833
$ws.Range("M126").Value = -8010.071599999999
$ws.Range("N126").Value = -12773
$ws.Range("H132").Value = 3243
$ws.Range("I132").Value = 2546
$ws.Range("K132").Value = 7638
$ws.Range("M132").Value = -5108
